# Auto-generated script to apply scheduled market-price refresh to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 372
$ws.Range("I5").Value = 287
$ws.Range("K5").Value = 287
$ws.Range("M5").Value = -172
# Row 17
$ws.Range("H17").Value = 561925
$ws.Range("I17").Value = 1999
$ws.Range("J17").Value = 617917.6
$ws.Range("K17").Value = 5997
$ws.Range("L17").Value = 1853752.8
$ws.Range("M17").Value = -5829
$ws.Range("N17").Value = -1854088.8
# Row 43
$ws.Range("H43").Value = 12606.777
$ws.Range("I43").Value = 1114.25
$ws.Range("J43").Value = 21800.8
$ws.Range("K43").Value = 1114.25
$ws.Range("L43").Value = 21800.8
$ws.Range("M43").Value = -1045.25
$ws.Range("N43").Value = -21938.8
# Row 82
$ws.Range("H82").Value = 1078.2858
$ws.Range("I82").Value = 1078.2858
$ws.Range("K82").Value = 3234.8574
$ws.Range("M82").Value = -2828.8574
# Row 85
$ws.Range("H85").Value = 1078.2858
$ws.Range("I85").Value = 1078.2858
$ws.Range("K85").Value = 3234.8574
$ws.Range("M85").Value = -1830.8574
# Row 99
$ws.Range("H99").Value = 2603.3635
$ws.Range("I99").Value = 402
$ws.Range("J99").Value = 6455.75
$ws.Range("K99").Value = 1206
$ws.Range("L99").Value = 19367.25
$ws.Range("M99").Value = 292
$ws.Range("N99").Value = -22363.25
# Row 132
$ws.Range("H132").Value = 68329
$ws.Range("I132").Value = 45124.78
$ws.Range("K132").Value = 135374.34
$ws.Range("M132").Value = -132844.34
# Row 137
$ws.Range("H137").Value = 2245.8096
$ws.Range("I137").Value = 2260.8
$ws.Range("J137").Value = 2208.3333
$ws.Range("K137").Value = 6782.400000000001
$ws.Range("L137").Value = 6624.999899999999
$ws.Range("M137").Value = -4232.400000000001
$ws.Range("N137").Value = -11724.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2636.56
$ws.Range("I32").Value = 2120.1738
$ws.Range("K32").Value = 2120.1738
$ws.Range("M32").Value = -1833.1738
# Row 110
$ws.Range("H110").Value = 15673.5
$ws.Range("I110").Value = 26825.143
$ws.Range("K110").Value = 26825.143
$ws.Range("M110").Value = -24780.143
# Row 122
$ws.Range("H122").Value = 2596.3333
$ws.Range("I122").Value = 2152.2
$ws.Range("K122").Value = 6456.599999999999
$ws.Range("M122").Value = -4006.599999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2166.6667
$ws.Range("I16").Value = 2000
$ws.Range("K16").Value = 2000
$ws.Range("M16").Value = -1713
# Row 31
$ws.Range("H31").Value = 50238.953
$ws.Range("I31").Value = 85287.336
$ws.Range("J31").Value = 8180.9
$ws.Range("K31").Value = 85287.336
$ws.Range("L31").Value = 8180.9
$ws.Range("M31").Value = -84992.336
$ws.Range("N31").Value = -8770.9
# Row 32
$ws.Range("H32").Value = 8555.333000000001
$ws.Range("J32").Value = 3379.6
$ws.Range("L32").Value = 3379.6
$ws.Range("N32").Value = -4011.6
# Row 34
$ws.Range("H34").Value = 50238.953
$ws.Range("I34").Value = 85287.336
$ws.Range("J34").Value = 8180.9
$ws.Range("K34").Value = 85287.336
$ws.Range("L34").Value = 8180.9
$ws.Range("M34").Value = -85085.336
$ws.Range("N34").Value = -8584.9
# Row 35
$ws.Range("H35").Value = 616.6667
$ws.Range("I35").Value = 540
$ws.Range("K35").Value = 540
$ws.Range("M35").Value = -246
# Row 36
$ws.Range("H36").Value = 16248.5
$ws.Range("J36").Value = 16664.666
$ws.Range("L36").Value = 16664.666
$ws.Range("N36").Value = -17440.666
# Row 40
$ws.Range("H40").Value = 16248.5
$ws.Range("J40").Value = 16664.666
$ws.Range("L40").Value = 16664.666
$ws.Range("N40").Value = -16984.666
# Row 58
$ws.Range("H58").Value = 2184.158
$ws.Range("I58").Value = 2319.375
$ws.Range("J58").Value = 1463
$ws.Range("K58").Value = 2319.375
$ws.Range("L58").Value = 1463
$ws.Range("M58").Value = -2116.375
$ws.Range("N58").Value = -1869
# Row 113
$ws.Range("H113").Value = 2166.6667
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
# Row 132
$ws.Range("H132").Value = 4282.1274
$ws.Range("I132").Value = 4223.6587
$ws.Range("J132").Value = 4681.6665
$ws.Range("K132").Value = 12670.9761
$ws.Range("L132").Value = 14044.9995
$ws.Range("M132").Value = -10140.9761
$ws.Range("N132").Value = -19104.9995
# Row 134
$ws.Range("H134").Value = 15685.807
$ws.Range("I134").Value = 7567.6924
$ws.Range("J134").Value = 57900
$ws.Range("K134").Value = 22703.0772
$ws.Range("L134").Value = 173700
$ws.Range("M134").Value = -20168.0772
$ws.Range("N134").Value = -178770
# Row 136
$ws.Range("H136").Value = 2184.158
$ws.Range("I136").Value = 2319.375
$ws.Range("J136").Value = 1463
$ws.Range("K136").Value = 6958.125
$ws.Range("L136").Value = 4389
$ws.Range("M136").Value = -4408.125
$ws.Range("N136").Value = -9489

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 955.7646999999999
$ws.Range("I26").Value = 157.53847
$ws.Range("K26").Value = 472.61541
$ws.Range("M26").Value = -184.61541
# Row 122
$ws.Range("H122").Value = 4150
$ws.Range("I122").Value = 3223.75
$ws.Range("J122").Value = 6002.5
$ws.Range("K122").Value = 29013.75
$ws.Range("L122").Value = 54022.5
$ws.Range("M122").Value = -26563.75
$ws.Range("N122").Value = -58922.5
# Row 136
$ws.Range("H136").Value = 920954.4399999999
$ws.Range("I136").Value = 3334499.8
$ws.Range("K136").Value = 10003499.4
$ws.Range("M136").Value = -9998399.399999999

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 21867.666
$ws.Range("I4").Value = 301.5
$ws.Range("K4").Value = 301.5
$ws.Range("M4").Value = -189.5
# Row 13
$ws.Range("H13").Value = 966.2
$ws.Range("I13").Value = 308.33334
$ws.Range("J13").Value = 1953
$ws.Range("K13").Value = 308.33334
$ws.Range("L13").Value = 1953
$ws.Range("M13").Value = -169.33334
$ws.Range("N13").Value = -2231
# Row 52
$ws.Range("H52").Value = 30000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 30000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -30518
# Row 132
$ws.Range("H132").Value = 315929.72
$ws.Range("I132").Value = 457848.03
$ws.Range("J132").Value = 3709.4
$ws.Range("K132").Value = 1373544.09
$ws.Range("L132").Value = 11128.2
$ws.Range("M132").Value = -1371014.09
$ws.Range("N132").Value = -16188.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 10627.154
$ws.Range("I7").Value = 12017.111
$ws.Range("K7").Value = 12017.111
$ws.Range("M7").Value = -11905.111
# Row 9
$ws.Range("H9").Value = 1658.4286
$ws.Range("I9").Value = 652.5
$ws.Range("K9").Value = 652.5
$ws.Range("M9").Value = -428.5
# Row 16
$ws.Range("H16").Value = 3885.7144
$ws.Range("I16").Value = 1139.2106
$ws.Range("J16").Value = 6154.5654
$ws.Range("K16").Value = 1139.2106
$ws.Range("L16").Value = 6154.5654
$ws.Range("M16").Value = -969.2106000000001
$ws.Range("N16").Value = -6494.5654
# Row 40
$ws.Range("H40").Value = 3964.4443
$ws.Range("I40").Value = 3334.6875
$ws.Range("K40").Value = 3334.6875
$ws.Range("M40").Value = -3198.6875
# Row 55
$ws.Range("H55").Value = 276.3
$ws.Range("J55").Value = 232
$ws.Range("L55").Value = 232
$ws.Range("N55").Value = -578
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
# Row 126
$ws.Range("H126").Value = 10627.154
$ws.Range("I126").Value = 12017.111
$ws.Range("K126").Value = 36051.333
$ws.Range("M126").Value = -33581.333
# Row 136
$ws.Range("H136").Value = 5774.2915
$ws.Range("I136").Value = 5307.3335
$ws.Range("K136").Value = 15922.0005
$ws.Range("M136").Value = -13372.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 33500
$ws.Range("J44").Value = 33500
$ws.Range("L44").Value = 33500
$ws.Range("N44").Value = -34608
# Row 50
$ws.Range("H50").Value = 12505452
$ws.Range("J50").Value = 12505452
$ws.Range("L50").Value = 12505452
$ws.Range("N50").Value = -12506714
# Row 96
$ws.Range("H96").Value = 6000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 6000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -8746
